# This survey workbook collects form responses on three tabs
# (enabling_environment, market_creation, organisational_change) and
# summarises them on the "all" tab. The edit appends the response rows that
# were captured once the submit buttons started wiring up responses into the
# sheet (see commit message: "add buttons functionality, download result,
# and other features"). market_creation and all already had (and keep) their
# existing content untouched.

$wb = $excel.ActiveWorkbook

# --- enabling_environment: three identical rows of (test) responses ---
$ee = $wb.Worksheets.Item("enabling_environment")
for ($r = 2; $r -le 4; $r++) {
    $ee.Range("A$r").Value = "Angleina Rianti"
    $ee.Range("B$r").Value = "a"
    $ee.Range("C$r").Value = "s"
    $ee.Range("D$r").Value = "s"
    $ee.Range("E$r").Value = "no"
    $ee.Range("F$r").Value = "no"
    $ee.Range("G$r").Value = "no"
    $ee.Range("R$r").Value = 0
}

# --- organisational_change: one real submitted response ---
$oc = $wb.Worksheets.Item("organisational_change")
$oc.Range("A2").Value = "Angelina Rianti"
$oc.Range("B2").Value = "angelinarianti23@gmail.com"
$oc.Range("C2").Value = "Swinburne"
$oc.Range("D2").Value = "School"
$oc.Range("E2").Value = "yes"
$oc.Range("F2").Value = "no"
$oc.Range("G2").Value = "yes"
$oc.Range("H2").Value = "Partially"
$oc.Range("I2").Value = "Entirely"
$oc.Range("J2").Value = "Entirely"
$oc.Range("K2").Value = "Partially"
$oc.Range("L2").Value = "Partially"
$oc.Range("Q2").Value = "Partially"
$oc.Range("R2").Value = "Partially"
$oc.Range("S2").Value = "Partially"
$oc.Range("T2").Value = 1
